$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.930.91"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.631.25"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Range("E3").Value = "  +3.81%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.23"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Range("E5").Value = "  +1.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.23"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.566"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.657.72"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Range("E9").Value = "  +4.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.23"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Range("E10").Value = "  +2.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.105"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Range("E11").Value = "  +2.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.336"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.098.37"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Range("E14").Value = "  +4.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.860.91"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.86"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000137"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.646.14"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Range("E18").Value = "  +4.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "346.92"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Range("E19").Value = "  +3.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.52"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.34"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Range("E21").Value = "  +2.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.18"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Range("E22").Value = "  +3.92%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.75"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.418"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Range("E25").Value = "  +2.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.993"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("E27").Value = "  +1.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0799"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.10"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Range("E29").Value = "  +2.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.997"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.28"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Range("E31").Value = "  +7.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.92"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Range("E32").Value = "  +2.06%  "
$ws.Range("E33").Value = "  +2.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.56"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.971"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Range("E35").Value = "  +5.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.00"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Range("E36").Value = "  +2.55%  "
$ws.Range("E37").Value = "  +2.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.72"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Range("E38").Value = "  +2.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.837"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.70"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Range("E40").Value = "  +5.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.41"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Range("E41").Value = "  +1.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0983"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.58"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Range("E46").Value = "  +5.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0529"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0229"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Range("E49").Value = "  +1.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.979.60"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Range("E50").Value = "  +4.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.66"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Range("E51").Value = "  +3.18%  "

# Rows 42-44: coin rotation (FirstDigitalUSD/Bittensor/Mantle cyclic shift) with updated data
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "277.48"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Range("E42").Value = "  -2.22%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.610"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Range("E43").Value = "  +1.82%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.994"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Range("E44").Value = "  -0.27%  "
